# Apply simulator full-month coverage / employee-name fix edit.
$wb = $excel.ActiveWorkbook

# --- "Weekly Timesheet" sheet ---
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Fix client names (B2:B6)
$ws1.Range("B2").Value = "McClure"
$ws1.Range("B3").Value = "Evans"
$ws1.Range("B4").Value = "Fritts"
$ws1.Range("B5").Value = "Hendricks"
$ws1.Range("B6").Value = "Regan"

# Populate Rate / Total for each day row
$ws1.Range("E2").Value = 95
$ws1.Range("F2").Value = 760
$ws1.Range("E3").Value = 95
$ws1.Range("F3").Value = 760
$ws1.Range("E4").Value = 95
$ws1.Range("F4").Value = 760
$ws1.Range("E5").Value = 95
$ws1.Range("F5").Value = 760
$ws1.Range("E6").Value = 95
$ws1.Range("F6").Value = 760

# Subtotal / Hourly subtotal / Grand total
$ws1.Range("F8").Value = 3800
$ws1.Range("F11").Value = 3800
$ws1.Range("F13").Value = 3800

# --- "Jason Schema" sheet ---
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Fix employee id (B2:B6)
$ws2.Range("B2").Value = "emp_4nlnrvy7"
$ws2.Range("B3").Value = "emp_4nlnrvy7"
$ws2.Range("B4").Value = "emp_4nlnrvy7"
$ws2.Range("B5").Value = "emp_4nlnrvy7"
$ws2.Range("B6").Value = "emp_4nlnrvy7"

# Fix client names (D2:D6)
$ws2.Range("D2").Value = "McClure"
$ws2.Range("D3").Value = "Evans"
$ws2.Range("D4").Value = "Fritts"
$ws2.Range("D5").Value = "Hendricks"
$ws2.Range("D6").Value = "Regan"

# Populate Rate / Total for each day row
$ws2.Range("F2").Value = 95
$ws2.Range("G2").Value = 760
$ws2.Range("F3").Value = 95
$ws2.Range("G3").Value = 760
$ws2.Range("F4").Value = 95
$ws2.Range("G4").Value = 760
$ws2.Range("F5").Value = 95
$ws2.Range("G5").Value = 760
$ws2.Range("F6").Value = 95
$ws2.Range("G6").Value = 760
